# Auto-generated PowerShell Excel COM-interop script
# Applies updated Betfair Back/Lay odds values for 2026-02-03 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 2.92
$ws.Range("G3").Value = 3.4
$ws.Range("I3").Value = 2.72
$ws.Range("J3").Value = 3.2
$ws.Range("N3").Value = 3.3
$ws.Range("P3").Value = 1.8
$ws.Range("S3").Value = 3.65
$ws.Range("V3").Value = 1.58
$ws.Range("W3").Value = 1.41
$ws.Range("Z3").Value = 20
$ws.Range("AN3").Value = 44
$ws.Range("F4").Value = 1.7
$ws.Range("G4").Value = 2.18
$ws.Range("H4").Value = 4.4
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 2.32
$ws.Range("K4").Value = 5.4
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 1.89
$ws.Range("O4").Value = 1.3
$ws.Range("P4").Value = 1.35
$ws.Range("Q4").Value = 2.16
$ws.Range("R4").Value = 1.08
$ws.Range("S4").Value = 2.92
$ws.Range("T4").Value = 2.02
$ws.Range("U4").Value = 1.43
$ws.Range("V4").Value = 1.13
$ws.Range("W4").Value = 1.85
$ws.Range("H5").Value = 7.6
$ws.Range("I5").Value = 13
$ws.Range("J5").Value = 3.8
$ws.Range("M5").Value = 1.08
$ws.Range("O5").Value = 1.23
$ws.Range("P5").Value = 1.59
$ws.Range("Q5").Value = 1.48
$ws.Range("R5").Value = 1.16
$ws.Range("S5").Value = 3.85
$ws.Range("X5").Value = 12
$ws.Range("F6").Value = 2.38
$ws.Range("N6").Value = 2.46
$ws.Range("T7").Value = 1.58
$ws.Range("Z8").Value = 16.5
$ws.Range("AD8").Value = 13
$ws.Range("AJ8").Value = 100
$ws.Range("AO8").Value = 10
$ws.Range("J9").Value = 3.55
$ws.Range("K9").Value = 3.9
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 3.9
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 1.99
$ws.Range("Q9").Value = 1.86
$ws.Range("R9").Value = 1.39
$ws.Range("S9").Value = 3.15
$ws.Range("T9").Value = 1.71
$ws.Range("U9").Value = 2.18
$ws.Range("V9").Value = 1.49
$ws.Range("W9").Value = 1.63
$ws.Range("X9").Value = 16.5
$ws.Range("Y9").Value = 13
$ws.Range("Z9").Value = 21
$ws.Range("AA9").Value = 48
$ws.Range("AB9").Value = 11.5
$ws.Range("AC9").Value = 8.4
$ws.Range("AE9").Value = 34
$ws.Range("AF9").Value = 16.5
$ws.Range("AH9").Value = 18
$ws.Range("AI9").Value = 50
$ws.Range("AK9").Value = 27
$ws.Range("AL9").Value = 48
$ws.Range("AM9").Value = 85
$ws.Range("AN9").Value = 20
$ws.Range("AO9").Value = 27
$ws.Range("H10").Value = 1.66
$ws.Range("I10").Value = 1.77
$ws.Range("J10").Value = 3.65
$ws.Range("K10").Value = 4
$ws.Range("L10").Value = 1.39
$ws.Range("N10").Value = 3.25
$ws.Range("O10").Value = 1.39
$ws.Range("Q10").Value = 2.14
$ws.Range("T10").Value = 2.02
$ws.Range("X10").Value = 14
$ws.Range("AB10").Value = 20
$ws.Range("AC10").Value = 9.800000000000001
$ws.Range("AJ10").Value = 220
$ws.Range("W11").Value = 1.3
$ws.Range("I12").Value = 4.9
$ws.Range("J12").Value = 3.5
$ws.Range("K12").Value = 3.85
$ws.Range("N12").Value = 3.65
$ws.Range("P12").Value = 1.92
$ws.Range("T12").Value = 1.79
$ws.Range("V12").Value = 1.26
$ws.Range("W12").Value = 1.87
$ws.Range("X12").Value = 17.5
$ws.Range("AC12").Value = 9.800000000000001
$ws.Range("AG12").Value = 1000
$ws.Range("AN12").Value = 17.5
$ws.Range("F13").Value = 1.86
$ws.Range("G13").Value = 1.99
$ws.Range("H13").Value = 3.95
$ws.Range("R13").Value = 1.55
$ws.Range("S13").Value = 2.48
$ws.Range("U13").Value = 2.34
$ws.Range("W13").Value = 2.02
$ws.Range("AA13").Value = 90
$ws.Range("AO13").Value = 40
$ws.Range("F14").Value = 2.58
$ws.Range("H14").Value = 2.74
$ws.Range("K14").Value = 3.65
$ws.Range("M14").Value = 1.07
$ws.Range("N14").Value = 3.75
$ws.Range("P14").Value = 1.95
$ws.Range("R14").Value = 1.34
$ws.Range("S14").Value = 3
$ws.Range("T14").Value = 1.73
$ws.Range("U14").Value = 2.2
$ws.Range("V14").Value = 1.5
$ws.Range("AE14").Value = 38
$ws.Range("AI14").Value = 50
$ws.Range("AL14").Value = 48
$ws.Range("AN14").Value = 28
$ws.Range("AO14").Value = 32
$ws.Range("F15").Value = 1.61
$ws.Range("G15").Value = 1.7
$ws.Range("H15").Value = 6.4
$ws.Range("J15").Value = 3.8
$ws.Range("K15").Value = 4.6
$ws.Range("L15").Value = 1.34
$ws.Range("N15").Value = 3.4
$ws.Range("O15").Value = 1.38
$ws.Range("P15").Value = 1.86
$ws.Range("Q15").Value = 1.98
$ws.Range("R15").Value = 1.31
$ws.Range("S15").Value = 3.5
$ws.Range("T15").Value = 1.98
$ws.Range("U15").Value = 1.83
$ws.Range("W15").Value = 2.44
$ws.Range("X15").Value = 16
$ws.Range("AA15").Value = 250
$ws.Range("AB15").Value = 8
$ws.Range("AD15").Value = 29
$ws.Range("AE15").Value = 130
$ws.Range("AG15").Value = 11.5
$ws.Range("AH15").Value = 28
$ws.Range("AI15").Value = 120
$ws.Range("AJ15").Value = 17.5
$ws.Range("AK15").Value = 21
$ws.Range("AL15").Value = 48
$ws.Range("AM15").Value = 180
$ws.Range("AN15").Value = 12.5
$ws.Range("AO15").Value = 190
$ws.Range("N16").Value = 4.1
$ws.Range("F17").Value = 1.86
$ws.Range("H17").Value = 4.2
$ws.Range("K17").Value = 4.2
$ws.Range("L17").Value = 1.3
$ws.Range("N17").Value = 3.95
$ws.Range("O17").Value = 1.27
$ws.Range("P17").Value = 2.12
$ws.Range("Q17").Value = 1.74
$ws.Range("R17").Value = 1.42
$ws.Range("S17").Value = 2.96
$ws.Range("T17").Value = 1.7
$ws.Range("U17").Value = 2.2
$ws.Range("X17").Value = 18
$ws.Range("Y17").Value = 24
$ws.Range("Z17").Value = 36
$ws.Range("AB17").Value = 10.5
$ws.Range("AC17").Value = 9.4
$ws.Range("AD17").Value = 17
$ws.Range("AE17").Value = 55
$ws.Range("AG17").Value = 10
$ws.Range("AH17").Value = 17.5
$ws.Range("AI17").Value = 60
$ws.Range("AK17").Value = 19
$ws.Range("AL17").Value = 34
$ws.Range("AM17").Value = 90
$ws.Range("AN17").Value = 12
$ws.Range("AO17").Value = 55
$ws.Range("F19").Value = 3.1
$ws.Range("G19").Value = 4.3
$ws.Range("H19").Value = 2.08
$ws.Range("I19").Value = 2.54
$ws.Range("N19").Value = 3
$ws.Range("O19").Value = 1.26
$ws.Range("P19").Value = 1.88
$ws.Range("Q19").Value = 1.9
$ws.Range("R19").Value = 1.26
$ws.Range("S19").Value = 2.28
$ws.Range("T19").Value = 1.06
$ws.Range("V19").Value = 1.65
$ws.Range("W19").Value = 1.31
$ws.Range("X19").Value = 21
$ws.Range("Z19").Value = 20
$ws.Range("AC19").Value = 11.5
$ws.Range("AE19").Value = 32
$ws.Range("AF19").Value = 34
$ws.Range("AJ19").Value = 90
$ws.Range("AK19").Value = 55
$ws.Range("F20").Value = 1.87
$ws.Range("G20").Value = 2.38
$ws.Range("H20").Value = 3.55
$ws.Range("I20").Value = 5.4
$ws.Range("J20").Value = 2.92
$ws.Range("N20").Value = 2.98
$ws.Range("P20").Value = 1.67
$ws.Range("R20").Value = 1.24
$ws.Range("T20").Value = 1.93
$ws.Range("W20").Value = 1.73
$ws.Range("X20").Value = 12.5
$ws.Range("Y20").Value = 19
$ws.Range("AF20").Value = 17.5
$ws.Range("AG20").Value = 15.5
$ws.Range("AJ20").Value = 38
$ws.Range("AK20").Value = 36
$ws.Range("AN20").Value = 30
$ws.Range("G21").Value = 5.6
$ws.Range("H21").Value = 1.81
$ws.Range("I21").Value = 1.91
$ws.Range("P21").Value = 1.85
$ws.Range("V21").Value = 2.08
$ws.Range("AL21").Value = 80
$ws.Range("AN21").Value = 100
$ws.Range("F22").Value = 1.89
$ws.Range("H22").Value = 5.3
$ws.Range("I22").Value = 6
$ws.Range("V22").Value = 1.2
$ws.Range("I23").Value = 5.9
$ws.Range("O23").Value = 1.6
$ws.Range("T23").Value = 2.36
$ws.Range("V23").Value = 1.2
$ws.Range("AJ23").Value = 29
$ws.Range("T24").Value = 2.3
$ws.Range("F25").Value = 3.4
$ws.Range("G25").Value = 5.7
$ws.Range("I25").Value = 2.5
$ws.Range("J25").Value = 2.74
$ws.Range("K25").Value = 4.8
$ws.Range("N25").Value = 2.32
$ws.Range("O25").Value = 1.23
$ws.Range("P25").Value = 1.47
$ws.Range("Q25").Value = 2.16
$ws.Range("S25").Value = 2.5
$ws.Range("V25").Value = 1.67
$ws.Range("W25").Value = 1.21
